$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data for this workbook lists football match rows where each
# row's "id" (column A) is a sequential position, but the actual match data
# (game id in column B through column AD) for some rows was mis-ordered.
# This update corrects the ordering by moving the match-data (columns B:AD)
# between rows, while column A (the sequential id) stays in place.
#
# Groups below are lists of row numbers describing a cycle: the data that
# was in the first row moves to the second row, the data in the second row
# moves to the third row, etc., and the data in the last row wraps back to
# the first row.

$groups = @(
    ,@(26, 27)
    ,@(38, 39)
    ,@(43, 44)
    ,@(61, 62)
    ,@(114, 115)
    ,@(119, 120, 121)
    ,@(123, 124)
    ,@(138, 139)
    ,@(156, 157)
    ,@(194, 195)
)

foreach ($rows in $groups) {
    # Snapshot the current (pre-edit) contents of columns B:AD for every row
    # in this cycle before writing anything back.
    $snapshots = @()
    foreach ($r in $rows) {
        $snapshots += ,($ws.Range("B$r`:AD$r").Value2)
    }

    $count = $rows.Count
    for ($i = 0; $i -lt $count; $i++) {
        $destRow = $rows[$i]
        $srcIndex = ($i - 1 + $count) % $count
        $ws.Range("B$destRow`:AD$destRow").Value2 = $snapshots[$srcIndex]
    }
}

Write-Host "row data reshuffle complete"
